# This script "rotates" the species-observation data among rows 4-13.
# Each target row receives the A/B/D/E/F/G/H/Q/R values that (before the
# edit) belonged to a different source row, per the mapping below.
# All other columns (C, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG,
# AT, AW, AX, AY, ...) are left untouched because they were already
# identical between source and target rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move together as a group.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# target row -> source row (values are copied from source's *original*
# state into target)
$mapping = @{
    4  = 7
    5  = 8
    6  = 13
    7  = 4
    8  = 5
    9  = 6
    10 = 9
    11 = 10
    12 = 11
    13 = 12
}

# Snapshot the original values for every row referenced above, before
# any writes happen (so chained overwrites don't corrupt later reads).
$snapshot = @{}
foreach ($row in $mapping.Values) {
    if (-not $snapshot.ContainsKey($row)) {
        $rowData = @{}
        foreach ($col in $cols) {
            $rowData[$col] = $ws.Range("$col$row").Value2
        }
        $snapshot[$row] = $rowData
    }
}

# Apply the snapshot values to their target rows.
foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $rowData = $snapshot[$sourceRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value = $rowData[$col]
    }
}
